$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the team record columns, styled like the rest of row 1
# (bold, bordered, centered — matching the existing header formatting).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Fill in the team's Wins/Losses/Ties for every player row (2-44).
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 94   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
